# "Edit alpha email format for gmail"
#
# The author pasted the email body in from Gmail, which flattens the
# original three paragraphs (+ the separate link-only paragraph) into a
# single Gmail-style paragraph that fakes blank lines with manual <w:br/>
# breaks, and re-styles every run as Arial / 12pt (sz=24 half-points) /
# color #222222, with a white run shading (Gmail's "quoted text" look).
# The CTA copy also changes from "~1 minute" to "~2 minutes", a new
# "-Jamie" signature paragraph is added, and a trailing empty paragraph
# is left at the end.
#
# Because nearly every run's formatting changes, the cleanest and most
# faithful way to reproduce this with COM is to build the exact target
# OOXML for the whole body and drop it in with Range.InsertXML, which
# replaces the contents of the range it's called on (here: the entire
# document body).

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Run properties shared by most of the body text runs.
$rPrPlain = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
# Same, plus the white shading Gmail adds behind visible text runs.
$rPrShaded = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'
# Bold "free" run (no shading, matches the diff).
$rPrBold = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="222222"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# Manual line break run (used in place of the blank-line paragraphs Gmail collapses).
$brRun = "<w:r>$rPrPlain<w:br/></w:r>"

# paragraph mark / default-run properties for the big paragraph.
$pPr1 = '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'
# paragraph mark / default-run properties for the blank + signature paragraphs.
$pPr2 = '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'

$para1Runs =
  "<w:r>$rPrShaded<w:t>Do you like trees, spies, games, or exercise?</w:t></w:r>" +
  $brRun + $brRun +
  "<w:r>$rPrShaded<w:t>If yes, I would like to introduce Sappy Secrets, a tree-spy-adventure game developed by me for you. You, as arborous-</w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r>$rPrShaded<w:t>espionager</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r>$rPrShaded<w:t xml:space=`"preserve`"> Agent Almond, solve mysteries and complete missions in a real-world romp through a </w:t></w:r>" +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  "<w:r>$rPrShaded<w:t>forest of lies. The game plays out wherever you like to walk or run, and you interact with it using your mobile device and earphones.</w:t></w:r>" +
  $brRun + $brRun +
  "<w:r>$rPrShaded<w:t>I need people to test the game as it is made. In exchange, you get to be part of the creation process and get the game for </w:t></w:r>" +
  "<w:r>$rPrBold<w:t>free</w:t></w:r>" +
  "<w:r>$rPrShaded<w:t>! Sign up here (~2 minutes):</w:t></w:r>" +
  $brRun + $brRun +
  "<w:r>$rPrShaded<w:t>https://forms.gle/DuaQHKz7fEjKnCRbA</w:t></w:r>"

$para1 = "<w:p $wNs>$pPr1$para1Runs</w:p>"
$para2 = "<w:p $wNs>$pPr2</w:p>"
$para3 = "<w:p $wNs>$pPr2<w:r>$rPrPlain<w:t>-Jamie</w:t></w:r></w:p>"
$para4 = "<w:p $wNs/>"

$xml = $para1 + $para2 + $para3 + $para4

# Replace the whole body's content in one shot so paragraph/run
# boundaries match the target exactly; Word keeps the trailing sectPr.
$d = $word.ActiveDocument
$d.Content.InsertXML($xml)
